$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A (country name) updates caused by re-sorting rows by updated case totals ---
$ws.Range("A1").Value = 'Datos actualizados a 22 de Mayo de 2020 a las 20:35'
$ws.Range("A73").Value = 'Sudan'
$ws.Range("A74").Value = 'Honduras'
$ws.Range("A96").Value = 'Mayotte'
$ws.Range("A97").Value = 'Nueva Zelanda'
$ws.Range("A98").Value = 'Eslovaquia'
$ws.Range("A114").Value = 'Zambia'
$ws.Range("A115").Value = 'Costa Rica'
$ws.Range("A116").Value = 'Venezuela'
$ws.Range("A153").Value = 'Mauritania'
$ws.Range("A154").Value = 'Birmania'
$ws.Range("A155").Value = 'Martinica'
$ws.Range("A156").Value = 'Islas Feroe'
$ws.Range("A173").Value = 'Malaui'
$ws.Range("A174").Value = 'Liechtenstein'
$ws.Range("A175").Value = 'Comoras'
$ws.Range("A176").Value = 'San Martin (Parte Holandesa)'
$ws.Range("A215").Value = 'San Bartolome'
$ws.Range("A216").Value = 'Bonaire, San Eustaquio y Saba'

# --- Updated statistics (Casos totales, Nuevos casos, Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes) ---
# Row 7
$ws.Range("B7").Value = 281904
$ws.Range("C7").Value = 1787
$ws.Range("E7").Value = 56318
$ws.Range("G7").Value = 688
$ws.Range("H7").Value = 28628
# Row 11
$ws.Range("B11").Value = 179584
$ws.Range("C11").Value = 563
$ws.Range("E11").Value = 12242
$ws.Range("G11").Value = 33
$ws.Range("H11").Value = 8342
# Row 17
$ws.Range("B17").Value = 82413
$ws.Range("C17").Value = 1089
$ws.Range("D17").Value = 42467
$ws.Range("E17").Value = 33701
$ws.Range("G17").Value = 93
$ws.Range("H17").Value = 6245
# Row 32
$ws.Range("B32").Value = 27892
$ws.Range("C32").Value = 994
$ws.Range("D32").Value = 13798
$ws.Range("E32").Value = 13853
$ws.Range("G32").Value = 4
$ws.Range("H32").Value = 241
# Row 33
$ws.Range("B33").Value = 24506
$ws.Range("C33").Value = 115
$ws.Range("E33").Value = 1854
$ws.Range("G33").Value = 9
$ws.Range("H33").Value = 1592
# Row 54
$ws.Range("B54").Value = 8414
$ws.Range("C54").Value = 240
$ws.Range("D54").Value = 4096
$ws.Range("E54").Value = 4306
# Row 73
$ws.Range("B73").Value = 3378
$ws.Range("C73").Value = 240
$ws.Range("D73").Value = 372
$ws.Range("E73").Value = 2869
$ws.Range("G73").Value = 16
$ws.Range("H73").Value = 137
# Row 74
$ws.Range("B74").Value = 3204
$ws.Range("C74").Value = 104
$ws.Range("D74").Value = 397
$ws.Range("E74").Value = 2651
$ws.Range("G74").Value = 5
$ws.Range("H74").Value = 156
# Row 96
$ws.Range("B96").Value = 1521
$ws.Range("C96").Value = 46
$ws.Range("D96").Value = 894
$ws.Range("E96").Value = 608
$ws.Range("H96").Value = 19
# Row 97
$ws.Range("B97").Value = 1504
$ws.Range("D97").Value = 1455
$ws.Range("E97").Value = 28
$ws.Range("H97").Value = 21
# Row 98
$ws.Range("B98").Value = 1503
$ws.Range("C98").Value = 1
$ws.Range("D98").Value = 1256
$ws.Range("E98").Value = 219
$ws.Range("H98").Value = 28
# Row 114
$ws.Range("B114").Value = 920
$ws.Range("C114").Value = 54
$ws.Range("D114").Value = 336
$ws.Range("E114").Value = 577
$ws.Range("H114").Value = 7
# Row 115
$ws.Range("B115").Value = 903
$ws.Range("D115").Value = 592
$ws.Range("E115").Value = 301
# Row 116
$ws.Range("B116").Value = 882
$ws.Range("D116").Value = 262
$ws.Range("E116").Value = 610
$ws.Range("H116").Value = 10
# Row 153
$ws.Range("B153").Value = 200
$ws.Range("C153").Value = 27
$ws.Range("D153").Value = 7
$ws.Range("E153").Value = 187
$ws.Range("G153").Value = 1
# Row 154
$ws.Range("B154").Value = 199
$ws.Range("D154").Value = 108
$ws.Range("E154").Value = 85
$ws.Range("H154").Value = 6
# Row 155
$ws.Range("B155").Value = 197
$ws.Range("C155").Value = 5
$ws.Range("D155").Value = 91
$ws.Range("E155").Value = 92
$ws.Range("H155").Value = 14
# Row 156
$ws.Range("B156").Value = 187
$ws.Range("D156").Value = 187
$ws.Range("E156").Value = 0
$ws.Range("H156").Value = 0
# Row 173
$ws.Range("C173").Value = 10
$ws.Range("D173").Value = 28
$ws.Range("E173").Value = 51
$ws.Range("H173").Value = 3
# Row 174
$ws.Range("B174").Value = 82
$ws.Range("C174").Value = 0
$ws.Range("D174").Value = 55
$ws.Range("E174").Value = 26
# Row 175
$ws.Range("B175").Value = 78
$ws.Range("C175").Value = 44
$ws.Range("D175").Value = 18
$ws.Range("E175").Value = 59
$ws.Range("H175").Value = 1
# Row 176
$ws.Range("B176").Value = 77
$ws.Range("D176").Value = 59
$ws.Range("E176").Value = 3
$ws.Range("H176").Value = 15
